$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J - copy formatting (style) from H1
# (xlPasteFormats = -4122) so the same cell style is reused, then set values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values per row for column I (I0) and column J (IF)
$data = @{
    2  = @(7, 7)
    3  = @(7, 7)
    4  = @(8, 9)
    5  = @(9, 9)
    6  = @(7, 8)
    7  = @(7, 8)
    8  = @(7, 8)
    9  = @(8, 8)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(8, 8)
    14 = @(6, 7)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(4, 4)
    19 = @(7, 7)
    20 = @(3, 3)
    21 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
